{"js": "// The underlying edit (per the diff) is a spelling correction in the\n// candidate's declaration paragraph: \"hearby\" -> \"hereby\". (Everything\n// else in the diff is Word re-saving the file and dropping/merging its\n// internal, invisible spell/grammar-check bookmarks - w:proofErr - and\n// the runs they split; the visible document text besides that single\n// word is identical before/after.)\n\nconst body = context.document.body;\n\n// Find the single, case-sensitive occurrence of the misspelling and fix it.\nconst misspelled = body.search(\"hearby\", { matchCase: true, matchWholeWord: true });\nmisspelled.load(\"text\");\nawait context.sync();\n\nif (misspelled.items.length > 0) {\n  for (let i = 0; i < misspelled.items.length; i++) {\n    misspelled.items[i].insertText(\"hereby\", Word.InsertLocation.replace);\n  }\n} else {\n  // Fallback in case matchWholeWord behaves unexpectedly in some hosts.\n  const fallback = body.search(\"hearby\", { matchCase: true });\n  fallback.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < fallback.items.length; i++) {\n    fallback.items[i].insertText(\"hereby\", Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The underlying edit (per the diff) is a spelling correction in the\n# candidate's declaration paragraph: \"hearby\" -> \"hereby\". (Everything\n# else in the diff is just Word re-saving the file and dropping/merging\n# its internal, invisible spell/grammar-check bookmarks - w:proofErr -\n# and the runs they split around; the visible document text besides\n# that single word is identical before and after.)\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"hearby\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n$find.Replacement.Text = \"hereby\"\n\n# Forward:=True, Wrap:=wdFindContinue(1), Format:=False, Replace:=wdReplaceAll(2)\n$find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
